$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "week_day" column (D) is being removed. Shift start_time/end_time
# one column to the left (D<-E, E<-F) instead of using Columns.Delete()
# so the existing merged <cols> width definition for columns 1-6 is left
# untouched (Columns.Delete() would otherwise shrink that range).
$ws.Range("D1").Value = $ws.Range("E1").Text
$ws.Range("D2").Value = $ws.Range("E2").Text
$ws.Range("E1").Value = $ws.Range("F1").Text
$ws.Range("E2").Value = ""

# New columns: start_date, num_weeks, repetition
$ws.Range("F1").Value = "start_date"
$ws.Range("F2").Value = "შეიყვანეთ yyyy-MM-dd ფორმატში"

$ws.Range("G1").Value = "num_weeks"
$ws.Range("G2").Value = "კვირების რაოდენობა"

$ws.Range("H1").Value = "repetition"
$ws.Range("H2").Value = "რამდენ კვირაში ერთხელ განმეორდეს (შეიყვანეთ რიცხვი)"

# Give the new "num_weeks" column its own width (closest achievable value
# to the authored 13.85546875 given this host's column-width rounding).
$ws.Columns.Item(7).ColumnWidth = 13
